# Auto-generated-by-assistant PowerShell COM-interop edit script
$wb = $excel.ActiveWorkbook

# ---- Sheet: تفاصيل المجموعات ----
$ws1 = $wb.Worksheets.Item("تفاصيل المجموعات")
$s1rows = @(
    @("المجموعة_1", "أصلية", 7, 210, 250, 160, 40000, 160),
    @("المجموعة_1", "أصلية", 5, 168, 235, 170, 39950, 182),
    @("المجموعة_2", "أصلية", 3, 145, 200, 77, 15400, 130),
    @("المجموعة_2", "أصلية", 1, 126, 170, 91, 15470, 182),
    @("المجموعة_2", "أصلية", 1, 126, 170, 91, 15470, 182),
    @("المجموعة_3", "بواقي عادية", 6, 200, 300, 1, 300, 12),
    @("المجموعة_3", "بواقي عادية", 6, 200, 300, 1, 300, 12),
    @("المجموعة_4", "بواقي عادية", 6, 200, 300, 1, 300, 10),
    @("المجموعة_4", "بواقي عادية", 6, 200, 300, 1, 300, 10),
    @("المجموعة_5", "بواقي عادية", 6, 200, 300, 1, 300, 8),
    @("المجموعة_5", "بواقي عادية", 6, 200, 300, 1, 300, 8),
    @("المجموعة_6", "بواقي عادية", 6, 200, 300, 1, 300, 6),
    @("المجموعة_6", "بواقي عادية", 6, 200, 300, 1, 300, 6),
    @("المجموعة_7", "بواقي عادية", 6, 200, 300, 1, 300, 4),
    @("المجموعة_7", "بواقي عادية", 6, 200, 300, 1, 300, 4),
    @("المجموعة_8", "بواقي عادية", 6, 200, 300, 1, 300, 2),
    @("المجموعة_8", "بواقي عادية", 6, 200, 300, 1, 300, 2),
    @("المجموعة_9", "بواقي عادية", 2, 133, 190, 1, 190, 12),
    @("المجموعة_9", "بواقي عادية", 2, 133, 190, 1, 190, 12),
    @("المجموعة_9", "بواقي عادية", 2, 133, 190, 1, 190, 12),
    @("المجموعة_10", "بواقي عادية", 2, 133, 190, 1, 190, 9),
    @("المجموعة_10", "بواقي عادية", 2, 133, 190, 1, 190, 9),
    @("المجموعة_10", "بواقي عادية", 2, 133, 190, 1, 190, 9),
    @("المجموعة_11", "بواقي عادية", 2, 133, 190, 1, 190, 6),
    @("المجموعة_11", "بواقي عادية", 2, 133, 190, 1, 190, 6),
    @("المجموعة_11", "بواقي عادية", 2, 133, 190, 1, 190, 6),
    @("المجموعة_12", "بواقي عادية", 2, 133, 190, 1, 190, 3),
    @("المجموعة_12", "بواقي عادية", 2, 133, 190, 1, 190, 3),
    @("المجموعة_12", "بواقي عادية", 2, 133, 190, 1, 190, 3)
)
$r = 2
foreach ($row in $s1rows) {
    $ws1.Cells.Item($r,1).Value = $row[0]
    $ws1.Cells.Item($r,2).Value = $row[1]
    $ws1.Cells.Item($r,3).Value = $row[2]
    $ws1.Cells.Item($r,4).Value = $row[3]
    $ws1.Cells.Item($r,5).Value = $row[4]
    $ws1.Cells.Item($r,6).Value = $row[5]
    $ws1.Cells.Item($r,7).Value = $row[6]
    $ws1.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}
$ws1.Cells.Item(31,1).Value = "المجموع"
$ws1.Cells.Item(31,3).Formula = "=SUM(C2:C30)"
$ws1.Cells.Item(31,4).Formula = "=SUM(D2:D30)"
$ws1.Cells.Item(31,5).Formula = "=SUM(E2:E30)"
$ws1.Cells.Item(31,6).Formula = "=SUM(F2:F30)"
$ws1.Cells.Item(31,7).Formula = "=SUM(G2:G30)"
$ws1.Cells.Item(31,8).Formula = "=SUM(H2:H30)"

# ---- Sheet: ملخص المجموعات ----
$ws2 = $wb.Worksheets.Item("ملخص المجموعات")
$s2rows = @(
    @("المجموعة_1", "أصلية", 378, 40000, 15111600, 2),
    @("المجموعة_2", "أصلية", 397, 15400, 6131440, 3),
    @("المجموعة_3", "بواقي عادية", 400, 300, 120000, 2),
    @("المجموعة_4", "بواقي عادية", 400, 300, 120000, 2),
    @("المجموعة_5", "بواقي عادية", 400, 300, 120000, 2),
    @("المجموعة_6", "بواقي عادية", 400, 300, 120000, 2),
    @("المجموعة_7", "بواقي عادية", 400, 300, 120000, 2),
    @("المجموعة_8", "بواقي عادية", 400, 300, 120000, 2),
    @("المجموعة_9", "بواقي عادية", 399, 190, 75810, 3),
    @("المجموعة_10", "بواقي عادية", 399, 190, 75810, 3),
    @("المجموعة_11", "بواقي عادية", 399, 190, 75810, 3),
    @("المجموعة_12", "بواقي عادية", 399, 190, 75810, 3)
)
$r = 2
foreach ($row in $s2rows) {
    $ws2.Cells.Item($r,1).Value = $row[0]
    $ws2.Cells.Item($r,2).Value = $row[1]
    $ws2.Cells.Item($r,3).Value = $row[2]
    $ws2.Cells.Item($r,4).Value = $row[3]
    $ws2.Cells.Item($r,5).Value = $row[4]
    $ws2.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}
$ws2.Cells.Item(14,1).Value = "المجموع"
$ws2.Cells.Item(14,3).Formula = "=SUM(C2:C13)"
$ws2.Cells.Item(14,4).Formula = "=SUM(D2:D13)"
$ws2.Cells.Item(14,5).Formula = "=SUM(E2:E13)"
$ws2.Cells.Item(14,6).Formula = "=SUM(F2:F13)"

# ---- Sheet: السجاد المتبقي ----
$ws3 = $wb.Worksheets.Item("السجاد المتبقي")
$ws3.Cells.Item(4,4).Value = 12

# ---- Sheet: ملخص الواجهة ----
$ws4 = $wb.Worksheets.Item("ملخص الواجهة")
$s4rows = @(
    @(2, 40000, 378, "المجموعة_1", "أصلية"),
    @(3, 15400, 397, "المجموعة_2", "أصلية"),
    @(2, 300, 400, "المجموعة_3", "بواقي عادية"),
    @(2, 300, 400, "المجموعة_4", "بواقي عادية"),
    @(2, 300, 400, "المجموعة_5", "بواقي عادية"),
    @(2, 300, 400, "المجموعة_6", "بواقي عادية"),
    @(2, 300, 400, "المجموعة_7", "بواقي عادية"),
    @(2, 300, 400, "المجموعة_8", "بواقي عادية"),
    @(3, 190, 399, "المجموعة_9", "بواقي عادية"),
    @(3, 190, 399, "المجموعة_10", "بواقي عادية"),
    @(3, 190, 399, "المجموعة_11", "بواقي عادية"),
    @(3, 190, 399, "المجموعة_12", "بواقي عادية")
)
$r = 2
foreach ($row in $s4rows) {
    $ws4.Cells.Item($r,1).Value = $row[0]
    $ws4.Cells.Item($r,2).Value = $row[1]
    $ws4.Cells.Item($r,3).Value = $row[2]
    $ws4.Cells.Item($r,4).Value = $row[3]
    $ws4.Cells.Item($r,5).Value = $row[4]
    $r = $r + 1
}
$ws4.Cells.Item(14,1).Formula = "=SUM(A2:A13)"
$ws4.Cells.Item(14,2).Formula = "=SUM(B2:B13)"
$ws4.Cells.Item(14,3).Formula = "=SUM(C2:C13)"

# ---- Sheet: الإجماليات ----
$ws5 = $wb.Worksheets.Item("الإجماليات")
$ws5.Cells.Item(2,2).Value = 2673160
$ws5.Cells.Item(2,3).Value = 24479080

# ---- Sheet: اقتراحات تشكيل مجموعات ----
$ws6 = $wb.Worksheets.Item("اقتراحات تشكيل مجموعات")
$ws6.Cells.Item(2,4).Value = 12
$ws6.Cells.Item(2,8).Value = "مرشح بطول 235 وبعدد ≈ 12"

# ---- Sheet: تدقيق الكميات ----
$ws7 = $wb.Worksheets.Item("تدقيق الكميات")
$ws7.Cells.Item(6,5).Value = 170
$ws7.Cells.Item(6,6).Value = 12

Write-Host "Done applying edits"
